$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing rows 5-7 (their data is being folded into rows 2-4 below)
$ws.Rows("5:7").Delete()

# New TPM-derived values for rows 2-4 (Sending cluster = FAPs, Ligand = Wnt1, Receptor = Fzd1)
$rowData = @(
    @(2, @("FAPs","Wnt1","Fzd1","ECs",3,1,0.3532066666666667,1.05962,1,1,3,1,0.8596446666666667,2.578934,0.05286426382906832,0.05286426382906832,0.3036322272311111,2.73269004508,0.05286426382906832,0.05286426382906832)),
    @(3, @("FAPs","Wnt1","Fzd1","FAPs",3,1,0.3532066666666667,1.05962,1,1,3,1,10.435983,31.307949,0.6417658132713033,0.6417658132713032,3.686058768819999,33.17452891937999,0.6417658132713033,0.6417658132713032)),
    @(4, @("FAPs","Wnt1","Fzd1","MuSCs",3,1,0.3532066666666667,1.05962,1,1,3,1,4.965729333333333,14.897188,0.3053699228996285,0.3053699228996284,1.753928705395555,15.78535834856,0.3053699228996285,0.3053699228996284))
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
